$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.285.17"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "2.941.23"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "569.39"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").Value = "159.29"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "2.937.47"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "34.52"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "65.316.97"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "3.425.05"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "7.02"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "2.937.81"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "15.46"
$ws.Range("E20").Value = "  +11.01%  "
$ws.Range("D21").Value = "445.36"
$ws.Range("E21").Value = "  -4.21%  "
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "82.45"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "12.17"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  -4.98%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "8.01"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").Value = "2.59"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").Value = "27.35"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "5.75"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "0.972"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "44.30"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  -8.98%  "
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "2.85"
$ws.Range("E42").Value = "  -7.63%  "
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").Value = "383.18"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "0.0353"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "2.697.83"
$ws.Range("E47").Value = "  -3.99%  "
$ws.Range("D48").Value = "134.01"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").Value = "  +5.04%  "
$ws.Range("D51").Value = "23.44"
$ws.Range("E51").Value = "  -0.71%  "
